$d = $word.ActiveDocument

# Change 1
$d.Content.Find.Execute("Entrepreneurship", $true, $false, $false, $false, $false, $true, 1, $false, "Entrepreneurship and Innovation", 2) | Out-Null

# Change 2
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2) | Out-Null

# Change 3
$d.Content.Find.Execute("Curso (semestre ideal): EF (5), EM (4), EA (2), EB (8), EP (6), EQD (5), EQN (8)", $true, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EF (5), EM (4), EA (2), EB (8), EP (6), EQN (8)", 2) | Out-Null

# Change 4
$d.Content.Find.Execute("Fomentar a cultura do empreendedorismo; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup ao longo do semestre.", $true, $false, $false, $false, $false, $true, 1, $false, "Fomentar a cultura do empreendedorismo e da Inovação; Desenvolver habilidades empreendedoras; Apresentar conhecimentos necessários para a criação de startups. A disciplina é aplicada através de Aprendizagem baseada em Projetos, onde o projeto a ser desenvolvido é da criação de uma startup com uma proposta de produto/serviço inovador ao longo do semestre.", 2) | Out-Null

# Change 5
$d.Content.Find.Execute("Promote the culture of entrepreneurship. Develop entrepreneurial skills. Present knowledge needed to create startups. The discipline is applied through Project-Based Learning, where the project to be developed is the creation of a startup during the semester.", $true, $false, $false, $false, $false, $true, 1, $false, "Foster a culture of entrepreneurship and innovation; Develop entrepreneurial skills; To present the necessary knowledge for the creation of startups. The discipline is applied through Project-Based Learning, where the project to be developed is the creation of a startup with an innovative product/service proposal throughout the semester.", 2) | Out-Null

# Change 6
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false, $true, 1, $false, "11079086 - Herlandí de Souza Andrade", 2) | Out-Null

# Change 7
$d.Content.Find.Execute("1.Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2.Modelo de Negócios (Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.3.Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente.4.Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Características do Comportamento Empreendedor: Busca de oportunidades e iniciativa. Correr riscos calculados. Exigência de qualidade e eficiência. Persistência. Comprometimento. Busca de informações. Estabelecimento de metas. Monitoramento e planejamento sistemático. Persuasão e rede contatos. Independência e autoconfiança.2. Estratégia, Inovação e Marketing.3. Design Thinking.4. Modelo de Negócios (Business Model Canvas e Lean Startup - Lean Canvas): Problema. Segmento de Clientes. Proposta de Valor Única. Solução. Métricas-Chave. Canais. Estrutura de Custos. Fluxos de Receita. Vantagem Injusta.5. Produto mínimo viável: Ciclo Construir-Mensurar-Aprender. Valor da vida útil do cliente. Prototipação rápida.6. Gestão de processos e Gerenciamento ágil de projetos.7. Plano de Negócios: Marketing, Finanças, Recursos Humanos, Desenvolvimento de Produtos e Tecnologia da Informação e Comunicação.8. Proposta da criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável e uma rodada de PITCH.9. Desenvolvimento de atividade prática extensionista (produção de conteúdo digital sobre empreendedorismo e inovação)10. Visita (viagem didática complementar) a um ambiente de inovação e empreendedorismo (ex. incubadora/aceleradora ou parque tecnológico), para compreender o desenvolvimento dos processos de empreendedorismo e inovação.", 2) | Out-Null

# Change 8
$d.Content.Find.Execute("1.Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Requirement of quality and efficiency. Persistence. Commitment. Search for information. Setting goals. Monitoring and systematic planning. Persuasion and network contacts. Independence and self-confidence.2.Business Model (Lean Canvas): Problem. Customer Segments. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage. 3.Minimum Viable Product: Build-Measure-Learn Cycle. Customer Lifetime Value.4.Business Plan: Marketing. Finance. Human Resources. Product Development. Information and communication technology.", $true, $false, $false, $false, $false, $true, 1, $false, "1. Characteristics of Entrepreneurial Behavior: Search for opportunities and initiative. Take calculated risks. Demand for quality and efficiency. Persistence. Commitment. Information search. Setting goals. Systematic monitoring and planning. Persuasion and networking. Independence and self-confidence.2. Strategy, Innovation and marketing.3. Design Thinking.4. Business Model (Business Model Canvas and Lean Startup - Lean Canvas): Problem. Customer Segment. Unique Value Proposition. Solution. Key Metrics. Channels. Cost Structure. Revenue Streams. Unfair Advantage.5. Minimum Viable Product: Build-Measure-Learn Cycle. Customer lifetime value. Rapid prototyping.6. Process Management and Agile Project Management7. Business Plan: Marketing, Finance, Human Resources, Product Development and Information and Communication Technology. 8. Proposal for the creation of a startup, from the business model to the business plan, including the assembly of the minimum viable product and a PITCH round. 9. Development of practical extension activity (production of digital content on entrepreneurship and innovation)10. Visit (complementary didactic trip) to an environment of innovation and entrepreneurship (eg incubator/accelerator or technology park), to understand the development of entrepreneurship and innovation processes.", 2) | Out-Null

# Change 9
$d.Content.Find.Execute("Os alunos montarão equipes que serão responsáveis por propor a criação de uma startup, do modelo de negócios ao plano de negócios, incluindo a montagem do produto mínimo viável.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.", 2) | Out-Null

# Change 10
$d.Content.Find.Execute("Avaliação dos trabalhos e apresentações ao longo do semestre", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas", 2) | Out-Null

# Change 11
$d.Content.Find.Execute("NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação", 2) | Out-Null
